$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Gehan Adel, Administrator'
$ws.Range("G3").Value = 'Dr. Veronia Rafat, Dr. Majorelle Magdy, Dr. Eman Tantawi, Administrator, Dr. Asmaa Reda, Dr. Hend Mahmoud'
$ws.Range("G4").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Majorelle Magdy, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Hend Mahmoud'
$ws.Range("G5").Value = 'Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Asmaa Reda, Dr. Eman Tantawi'
$ws.Range("G6").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Majorelle Magdy, Dr. Mohammad El-Tanany, Dr. Manar Montaser, Dr. Alshimaa Atef'
$ws.Range("G7").Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Kerelos Zareef, Dr. Nada Mohammad, Dr. Abeer Ragab, Dr. Lamiaa Ossama, Dr. Amera Ahmad Saad'
$ws.Range("G8").Value = 'Dr. Abeer Ragab, Dr. Nada Mohammad'
$ws.Range("G9").Value = 'Dr. Shimaa Ashraf, Dr. Safa Hany'
$ws.Range("G12").Value = 'Dr. Madeha Saeed, Dr. Dina Adel, Dr. Amira Ibrahim, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Marina Youhanna'
$ws.Range("G13").Value = 'Dr. Amira Ibrahim, Dr. Yasmeena Fattoh, Dr. Esraa Mostafa'
$ws.Range("G27").Value = 'Dr. Hana Amr, Dr. Nourham Mostafa'
$ws.Range("G28").Value = 'Dr. Maryam Ashraf, Dr. Aya Emad'
